# Fix 2-digit years to 4-digit years in the "Creation date" column (M)
# e.g. "8/3/18 19:27:00" -> "8/3/2018 19:27:00"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 186; $r++) {
    $cell = $ws.Cells.Item($r, 13)
    $v = $cell.Value()

    if ($v -eq $null) { continue }

    $spaceIdx = $v.IndexOf(' ')
    if ($spaceIdx -lt 0) { continue }

    $datePart = $v.Substring(0, $spaceIdx)
    $timePart = $v.Substring($spaceIdx + 1)

    $dateParts = $datePart -split '/'
    if ($dateParts.Length -ne 3) { continue }

    $year = $dateParts[2]
    if ($year.Length -eq 2) {
        $newYear = "20$year"
        $newDate = "$($dateParts[0])/$($dateParts[1])/$newYear $timePart"
        $cell.Value = $newDate
    }
}
